{"js": "// Replace the twenty-five \"two-digit \u00d7 two-digit =\" expressions in the\n// document (the date line at the top is left untouched) with the new\n// values from the commit. Each old expression is unique in the document,\n// so a plain text search-and-replace is unambiguous for every pair.\nconst replacements = [\n  [\"23\u00d732=\", \"70\u00d711=\"],\n  [\"19\u00d767=\", \"33\u00d745=\"],\n  [\"99\u00d735=\", \"39\u00d755=\"],\n  [\"62\u00d737=\", \"83\u00d721=\"],\n  [\"25\u00d732=\", \"85\u00d776=\"],\n  [\"91\u00d712=\", \"75\u00d718=\"],\n  [\"57\u00d788=\", \"65\u00d753=\"],\n  [\"86\u00d772=\", \"86\u00d778=\"],\n  [\"95\u00d787=\", \"39\u00d731=\"],\n  [\"53\u00d720=\", \"57\u00d769=\"],\n  [\"93\u00d782=\", \"90\u00d757=\"],\n  [\"68\u00d726=\", \"65\u00d735=\"],\n  [\"48\u00d765=\", \"86\u00d753=\"],\n  [\"98\u00d720=\", \"82\u00d777=\"],\n  [\"35\u00d747=\", \"73\u00d741=\"],\n  [\"55\u00d716=\", \"37\u00d729=\"],\n  [\"66\u00d761=\", \"41\u00d793=\"],\n  [\"86\u00d782=\", \"66\u00d795=\"],\n  [\"74\u00d742=\", \"98\u00d794=\"],\n  [\"59\u00d736=\", \"62\u00d723=\"],\n  [\"64\u00d786=\", \"19\u00d782=\"],\n  [\"20\u00d755=\", \"93\u00d776=\"],\n  [\"26\u00d773=\", \"78\u00d771=\"],\n  [\"95\u00d786=\", \"55\u00d754=\"],\n  [\"60\u00d751=\", \"93\u00d783=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the twenty-five \"two-digit \u00d7 two-digit =\" expressions in the\n# document (the date line at the top is left untouched) with the new\n# values from the commit. Each old expression is unique in the document,\n# so Find/Replace with MatchWholeWord off but exact text match is\n# unambiguous for every pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"23\u00d732=\", \"70\u00d711=\"),\n    @(\"19\u00d767=\", \"33\u00d745=\"),\n    @(\"99\u00d735=\", \"39\u00d755=\"),\n    @(\"62\u00d737=\", \"83\u00d721=\"),\n    @(\"25\u00d732=\", \"85\u00d776=\"),\n    @(\"91\u00d712=\", \"75\u00d718=\"),\n    @(\"57\u00d788=\", \"65\u00d753=\"),\n    @(\"86\u00d772=\", \"86\u00d778=\"),\n    @(\"95\u00d787=\", \"39\u00d731=\"),\n    @(\"53\u00d720=\", \"57\u00d769=\"),\n    @(\"93\u00d782=\", \"90\u00d757=\"),\n    @(\"68\u00d726=\", \"65\u00d735=\"),\n    @(\"48\u00d765=\", \"86\u00d753=\"),\n    @(\"98\u00d720=\", \"82\u00d777=\"),\n    @(\"35\u00d747=\", \"73\u00d741=\"),\n    @(\"55\u00d716=\", \"37\u00d729=\"),\n    @(\"66\u00d761=\", \"41\u00d793=\"),\n    @(\"86\u00d782=\", \"66\u00d795=\"),\n    @(\"74\u00d742=\", \"98\u00d794=\"),\n    @(\"59\u00d736=\", \"62\u00d723=\"),\n    @(\"64\u00d786=\", \"19\u00d782=\"),\n    @(\"20\u00d755=\", \"93\u00d776=\"),\n    @(\"26\u00d773=\", \"78\u00d771=\"),\n    @(\"95\u00d786=\", \"55\u00d754=\"),\n    @(\"60\u00d751=\", \"93\u00d783=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
